# feat: add 2022-Q1 data
#
# Plan:
#  1. Duplicate the existing "总计" sheet (Worksheet.Copy) so we have a perfect
#     style-for-style clone placed right after it; that clone will become the
#     NEW "总计" sheet (with the 2022-Q1 row added on top).
#  2. Rename the ORIGINAL "总计" sheet to "2022-Q1" and turn it into the new
#     per-fund holdings sheet (matching the layout already used by the other
#     quarterly sheets, e.g. "2021-Q4").
#  3. Rename the cloned sheet back to "总计" and rewrite its table to include
#     the new 2022-Q1 summary row at the top, shifting the rest down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# --- Step 1: clone "总计" (placed immediately after the original) ---------
$totalSheet.Copy($null, $totalSheet)
$totalClone = $wb.Worksheets.Item($totalSheet.Index + 1)

# ===========================================================================
# Step 2: turn the original "总计" sheet into the "2022-Q1" fund-holdings sheet
# ===========================================================================
$q1 = $totalSheet
$q1.Name = "2022-Q1"

# Pull the header row (and its formatting) from an existing quarterly sheet so
# the new columns E:H show up with the same style (s=2) as B:D already have.
# (Column A of the header row is intentionally blank in every quarterly sheet,
# so only copy B1:H1.)
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy($q1.Range("B1:H1"))

# Drop the leftover rows (previously 总计's 5 data rows); only 2 funds this
# quarter.
$q1.Rows("4:6").Clear()

# Force text storage for the numeric-looking columns (fund code + the
# percentage-ish figures) so leading zeros / exact text are preserved, same
# as every other quarterly sheet.
$q1.Range("B2:B3").NumberFormat = "@"
$q1.Range("D2:G3").NumberFormat = "@"

# Row 2: 009956 / 广发恒誉混合A
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "009956"
$q1.Cells.Item(2, 3).Value = "广发恒誉混合A"
$q1.Cells.Item(2, 4).Value = "4.94"
$q1.Cells.Item(2, 5).Value = "21.40"
$q1.Cells.Item(2, 6).Value = "0.91"
$q1.Cells.Item(2, 7).Value = "0.0450"
$q1.Cells.Item(2, 8).Value = 7

# Row 3: 009957 / 广发恒誉混合C
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "009957"
$q1.Cells.Item(3, 3).Value = "广发恒誉混合C"
$q1.Cells.Item(3, 4).Value = "0.10"
$q1.Cells.Item(3, 5).Value = "21.40"
$q1.Cells.Item(3, 6).Value = "0.91"
$q1.Cells.Item(3, 7).Value = "0.0009"
$q1.Cells.Item(3, 8).Value = 7

# ===========================================================================
# Step 3: rebuild the NEW "总计" sheet (the clone) with the 2022-Q1 row added
# ===========================================================================
$total = $totalClone
$total.Name = "总计"

# Extend the styled index column (A) by one row, reusing A6's style (s=2) so
# the new row 7 matches the rest instead of minting a new style.
$total.Range("A6").Copy($total.Range("A7"))

$dates = @("2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$counts = @(2, 10, 2, 12, 10, 12)
$values = @(0.05, 2.42, 0.71, 2.21, 1.52, 3.9)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $dates[$i]
    $total.Cells.Item($r, 3).Value = $counts[$i]
    $total.Cells.Item($r, 4).Value = $values[$i]
}

$wb.Worksheets.Item("2021-Q4").Select()
